$wb = $excel.ActiveWorkbook

# 1. Remove the empty 'Planilha1' sheet (query refresh dropped the old preview sheet)
$planilha = $wb.Worksheets.Item("Planilha1")
$planilha.Delete()

$ws = $wb.Worksheets.Item("Chess_Blunders_and_Time")

# 2. Rename the external-data defined name (query refresh incremented the suffix)
$dname = $wb.Names.Item("Chess_Blunders_and_Time!DadosExternos_1")
$dname.Name = "DadosExternos_2"

# 3. Rename the table to match the refreshed query result name
$tbl = $ws.ListObjects.Item(1)
$tbl.Name = "part_00000_c47cebe9_98bf_4e9c_8a41_548b4e4c1333_c000"

# 4. Write the refreshed query data (same rows, new order/values per the latest pull)
$colA = New-Object 'object[,]' 90,1
$colDEF = New-Object 'object[,]' 90,3
$colA[0,0] = 45017
$colDEF[0,0] = 15633037
$colDEF[0,1] = 7327811
$colDEF[0,2] = 0.4687387997610445
$colA[1,0] = 45017
$colDEF[1,0] = 10372452
$colDEF[1,1] = 3054934
$colDEF[1,2] = 0.2945238020865269
$colA[2,0] = 45017
$colDEF[2,0] = 231753
$colDEF[2,1] = 261550
$colDEF[2,2] = 1.1285722299171963
$colA[3,0] = 45017
$colDEF[3,0] = 6100152
$colDEF[3,1] = 4825720
$colDEF[3,2] = 0.7910819271388647
$colA[4,0] = 45017
$colDEF[4,0] = 419130
$colDEF[4,1] = 247418
$colDEF[4,2] = 0.5903132679598215
$colA[5,0] = 45017
$colDEF[5,0] = 6121872
$colDEF[5,1] = 1954313
$colDEF[5,2] = 0.3192345413298416
$colA[6,0] = 45017
$colDEF[6,0] = 6767918
$colDEF[6,1] = 624557
$colDEF[6,2] = 0.09228199868851839
$colA[7,0] = 45017
$colDEF[7,0] = 19931
$colDEF[7,1] = 19589
$colDEF[7,2] = 0.982840800762631
$colA[8,0] = 45017
$colDEF[8,0] = 894097
$colDEF[8,1] = 585744
$colDEF[8,2] = 0.6551235492345909
$colA[9,0] = 45017
$colDEF[9,0] = 53893
$colDEF[9,1] = 9635
$colDEF[9,2] = 0.17878017553300057
$colA[10,0] = 45017
$colDEF[10,0] = 26587583
$colDEF[10,1] = 6660859
$colDEF[10,2] = 0.2505251793666239
$colA[11,0] = 45017
$colDEF[11,0] = 19306856
$colDEF[11,1] = 2796800
$colDEF[11,2] = 0.1448604578601508
$colA[12,0] = 45017
$colDEF[12,0] = 445505
$colDEF[12,1] = 430878
$colDEF[12,2] = 0.9671675963232735
$colA[13,0] = 45017
$colDEF[13,0] = 7886687
$colDEF[13,1] = 4483021
$colDEF[13,2] = 0.5684289233235704
$colA[14,0] = 45017
$colDEF[14,0] = 595179
$colDEF[14,1] = 189487
$colDEF[14,2] = 0.31836976775054227
$colA[15,0] = 44958
$colDEF[15,0] = 15581007
$colDEF[15,1] = 8142843
$colDEF[15,2] = 0.5226133971957011
$colA[16,0] = 44958
$colDEF[16,0] = 10355349
$colDEF[16,1] = 3258507
$colDEF[16,2] = 0.31466896963105734
$colA[17,0] = 44958
$colDEF[17,0] = 260646
$colDEF[17,1] = 321669
$colDEF[17,2] = 1.2341221426762736
$colA[18,0] = 44958
$colDEF[18,0] = 6229747
$colDEF[18,1] = 5471946
$colDEF[18,2] = 0.878357660431475
$colA[19,0] = 44958
$colDEF[19,0] = 407192
$colDEF[19,1] = 246947
$colDEF[19,2] = 0.6064632900449911
$colA[20,0] = 44958
$colDEF[20,0] = 5817154
$colDEF[20,1] = 1929773
$colDEF[20,2] = 0.33173833802577685
$colA[21,0] = 44958
$colDEF[21,0] = 6053959
$colDEF[21,1] = 568479
$colDEF[21,2] = 0.09390202345275216
$colA[22,0] = 44958
$colDEF[22,0] = 18802
$colDEF[22,1] = 19970
$colDEF[22,2] = 1.0621210509520265
$colA[23,0] = 44958
$colDEF[23,0] = 869950
$colDEF[23,1] = 577089
$colDEF[23,2] = 0.6633588137249267
$colA[24,0] = 44958
$colDEF[24,0] = 48769
$colDEF[24,1] = 11299
$colDEF[24,2] = 0.23168406159650598
$colA[25,0] = 44958
$colDEF[25,0] = 25557690
$colDEF[25,1] = 6724886
$colDEF[25,2] = 0.2631257363243705
$colA[26,0] = 44958
$colDEF[26,0] = 18223321
$colDEF[26,1] = 2668661
$colDEF[26,2] = 0.14644207825785432
$colA[27,0] = 44958
$colDEF[27,0] = 445972
$colDEF[27,1] = 445466
$colDEF[27,2] = 0.9988653996215009
$colA[28,0] = 44958
$colDEF[28,0] = 7696600
$colDEF[28,1] = 4566808
$colDEF[28,2] = 0.5933539484967388
$colA[29,0] = 44958
$colDEF[29,0] = 583051
$colDEF[29,1] = 190160
$colDEF[29,2] = 0.3261464262989001
$colA[30,0] = 44927
$colDEF[30,0] = 15787311
$colDEF[30,1] = 8071953
$colDEF[30,2] = 0.5112937219010888
$colA[31,0] = 44927
$colDEF[31,0] = 10289881
$colDEF[31,1] = 3154755
$colDEF[31,2] = 0.3065880936815499
$colA[32,0] = 44927
$colDEF[32,0] = 263151
$colDEF[32,1] = 333166
$colDEF[32,2] = 1.266063970876037
$colA[33,0] = 44927
$colDEF[33,0] = 6027622
$colDEF[33,1] = 5237607
$colDEF[33,2] = 0.8689342165119179
$colA[34,0] = 44927
$colDEF[34,0] = 420769
$colDEF[34,1] = 237509
$colDEF[34,2] = 0.5644641121375387
$colA[35,0] = 44927
$colDEF[35,0] = 6473761
$colDEF[35,1] = 2120866
$colDEF[35,2] = 0.32760956111910833
$colA[36,0] = 44927
$colDEF[36,0] = 6595429
$colDEF[36,1] = 621297
$colDEF[36,2] = 0.09420115052409783
$colA[37,0] = 44927
$colDEF[37,0] = 19816
$colDEF[37,1] = 20452
$colDEF[37,2] = 1.0320952765442066
$colA[38,0] = 44927
$colDEF[38,0] = 938294
$colDEF[38,1] = 617570
$colDEF[38,2] = 0.6581838954528112
$colA[39,0] = 44927
$colDEF[39,0] = 43913
$colDEF[39,1] = 9422
$colDEF[39,2] = 0.21456060847584998
$colA[40,0] = 44927
$colDEF[40,0] = 27433796
$colDEF[40,1] = 7147393
$colDEF[40,2] = 0.2605324104618989
$colA[41,0] = 44927
$colDEF[41,0] = 19437298
$colDEF[41,1] = 2867190
$colDEF[41,2] = 0.14750970016511555
$colA[42,0] = 44927
$colDEF[42,0] = 459455
$colDEF[42,1] = 456866
$colDEF[42,2] = 0.9943650629550228
$colA[43,0] = 44927
$colDEF[43,0] = 8052188
$colDEF[43,1] = 4737038
$colDEF[43,2] = 0.5882920269621127
$colA[44,0] = 44927
$colDEF[44,0] = 616624
$colDEF[44,1] = 191602
$colDEF[44,2] = 0.3107274449259192
$colA[45,0] = 45078
$colDEF[45,0] = 14181616
$colDEF[45,1] = 6157363
$colDEF[45,2] = 0.4341792218883941
$colA[46,0] = 45078
$colDEF[46,0] = 9227153
$colDEF[46,1] = 2513713
$colDEF[46,2] = 0.27242563334541003
$colA[47,0] = 45078
$colDEF[47,0] = 195757
$colDEF[47,1] = 205627
$colDEF[47,2] = 1.0504196529370597
$colA[48,0] = 45078
$colDEF[48,0] = 5467069
$colDEF[48,1] = 3921425
$colDEF[48,2] = 0.7172810513275029
$colA[49,0] = 45078
$colDEF[49,0] = 372902
$colDEF[49,1] = 200927
$colDEF[49,2] = 0.5388198507918971
$colA[50,0] = 45078
$colDEF[50,0] = 6039197
$colDEF[50,1] = 1994416
$colDEF[50,2] = 0.33024522962241504
$colA[51,0] = 45078
$colDEF[51,0] = 6972028
$colDEF[51,1] = 634464
$colDEF[51,2] = 0.09100135570310389
$colA[52,0] = 45078
$colDEF[52,0] = 19511
$colDEF[52,1] = 19502
$colDEF[52,2] = 0.999538721746707
$colA[53,0] = 45078
$colDEF[53,0] = 882404
$colDEF[53,1] = 566198
$colDEF[53,2] = 0.6416539362922199
$colA[54,0] = 45078
$colDEF[54,0] = 50718
$colDEF[54,1] = 10404
$colDEF[54,2] = 0.20513427185614574
$colA[55,0] = 45078
$colDEF[55,0] = 25947275
$colDEF[55,1] = 6285424
$colDEF[55,2] = 0.24223830826165754
$colA[56,0] = 45078
$colDEF[56,0] = 18562207
$colDEF[56,1] = 2626585
$colDEF[56,2] = 0.1415017621557609
$colA[57,0] = 45078
$colDEF[57,0] = 409992
$colDEF[57,1] = 375638
$colDEF[57,2] = 0.9162081211340709
$colA[58,0] = 45078
$colDEF[58,0] = 7504169
$colDEF[58,1] = 4010100
$colDEF[58,2] = 0.534382954328454
$colA[59,0] = 45078
$colDEF[59,0] = 597572
$colDEF[59,1] = 185971
$colDEF[59,2] = 0.31121103398418937
$colA[60,0] = 44986
$colDEF[60,0] = 17162966
$colDEF[60,1] = 8366104
$colDEF[60,2] = 0.48745094525037225
$colA[61,0] = 44986
$colDEF[61,0] = 11444714
$colDEF[61,1] = 3450796
$colDEF[61,2] = 0.30151876228623975
$colA[62,0] = 44986
$colDEF[62,0] = 267200
$colDEF[62,1] = 304133
$colDEF[62,2] = 1.1382223053892215
$colA[63,0] = 44986
$colDEF[63,0] = 6707749
$colDEF[63,1] = 5447409
$colDEF[63,2] = 0.8121068632711212
$colA[64,0] = 44986
$colDEF[64,0] = 451927
$colDEF[64,1] = 266027
$colDEF[64,2] = 0.5886503793754301
$colA[65,0] = 44986
$colDEF[65,0] = 6341716
$colDEF[65,1] = 2067588
$colDEF[65,2] = 0.3260297370617038
$colA[66,0] = 44986
$colDEF[66,0] = 6878413
$colDEF[66,1] = 624876
$colDEF[66,2] = 0.0908459553097495
$colA[67,0] = 44986
$colDEF[67,0] = 20206
$colDEF[67,1] = 21165
$colDEF[67,2] = 1.0474611501534197
$colA[68,0] = 44986
$colDEF[68,0] = 932382
$colDEF[68,1] = 601219
$colDEF[68,2] = 0.6448204705796551
$colA[69,0] = 44986
$colDEF[69,0] = 60661
$colDEF[69,1] = 10739
$colDEF[69,2] = 0.17703301956776182
$colA[70,0] = 44986
$colDEF[70,0] = 27936956
$colDEF[70,1] = 6934707
$colDEF[70,2] = 0.2482270079818288
$colA[71,0] = 44986
$colDEF[71,0] = 20184543
$colDEF[71,1] = 2900063
$colDEF[71,2] = 0.14367741692244407
$colA[72,0] = 44986
$colDEF[72,0] = 480277
$colDEF[72,1] = 461790
$colDEF[72,2] = 0.9615076299718704
$colA[73,0] = 44986
$colDEF[73,0] = 8327834
$colDEF[73,1] = 4672687
$colDEF[73,2] = 0.5610927163053442
$colA[74,0] = 44986
$colDEF[74,0] = 664570
$colDEF[74,1] = 192897
$colDEF[74,2] = 0.2902583625502204
$colA[75,0] = 45047
$colDEF[75,0] = 15684393
$colDEF[75,1] = 7408625
$colDEF[75,2] = 0.472356501140975
$colA[76,0] = 45047
$colDEF[76,0] = 10420495
$colDEF[76,1] = 3000185
$colDEF[76,2] = 0.2879119466013851
$colA[77,0] = 45047
$colDEF[77,0] = 225016
$colDEF[77,1] = 251946
$colDEF[77,2] = 1.1196803782842109
$colA[78,0] = 45047
$colDEF[78,0] = 6058432
$colDEF[78,1] = 4734187
$colDEF[78,2] = 0.7814211664008113
$colA[79,0] = 45047
$colDEF[79,0] = 391919
$colDEF[79,1] = 224057
$colDEF[79,2] = 0.571692109849229
$colA[80,0] = 45047
$colDEF[80,0] = 6346320
$colDEF[80,1] = 2115637
$colDEF[80,2] = 0.33336437494484994
$colA[81,0] = 45047
$colDEF[81,0] = 7203927
$colDEF[81,1] = 670283
$colDEF[81,2] = 0.09304411330098153
$colA[82,0] = 45047
$colDEF[82,0] = 20739
$colDEF[82,1] = 20948
$colDEF[82,2] = 1.0100776315155022
$colA[83,0] = 45047
$colDEF[83,0] = 934510
$colDEF[83,1] = 614671
$colDEF[83,2] = 0.6577468405902559
$colA[84,0] = 45047
$colDEF[84,0] = 53232
$colDEF[84,1] = 11849
$colDEF[84,2] = 0.2225916741809438
$colA[85,0] = 45047
$colDEF[85,0] = 27530595
$colDEF[85,1] = 7048181
$colDEF[85,2] = 0.25601266518213645
$colA[86,0] = 45047
$colDEF[86,0] = 19971533
$colDEF[86,1] = 2925932
$colDEF[86,2] = 0.1465051280740442
$colA[87,0] = 45047
$colDEF[87,0] = 443137
$colDEF[87,1] = 430167
$colDEF[87,2] = 0.9707313990932827
$colA[88,0] = 45047
$colDEF[88,0] = 8036388
$colDEF[88,1] = 4605267
$colDEF[88,2] = 0.5730518486663412
$colA[89,0] = 45047
$colDEF[89,0] = 613978
$colDEF[89,1] = 203448
$colDEF[89,2] = 0.3313604070504155

$ws.Range("A2:A91").Value = $colA
$ws.Range("D2:F91").Value = $colDEF

# 5. Restore the active selection noted after the refresh
$ws.Range("J26").Select()

Write-Host 'Edit complete'